$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D12").Value = "2016-03-08 06:22:43"
$wsZhCn.Range("G12").Value = "2016-03-08 06:23:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D12").Value = "2016-03-08 06:22:45"
$wsDeDe.Range("G12").Value = "2016-03-08 06:23:08"
